$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 (shifts old rows 8-11 down to 9-12),
# carrying the formatting of the row above (row 7, the "Ref" row).
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(8).Insert()

# Overwrite the new row 8 with the "Force" flag row content.
$ws.Range("A8").Value = "Force"
$ws.Range("B8:J8").Value = $false

Write-Output "done"
